$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new "Custom_scoring" column (D) ---
$ws.Range("D1").Value = "Custom_scoring"

# --- Row 2: text round with new valid-answer text "TON" and Custom_scoring 4 ---
$ws.Range("B2").Value = "text"
$ws.Range("C2").Value = "TON"
$ws.Range("D2").Value = 4

# --- Row 3: text round with new valid-answer text "QUEUE" and Custom_scoring 6 ---
$ws.Range("B3").Value = "text"
$ws.Range("C3").Value = "QUEUE"
$ws.Range("D3").Value = 6

# --- Rows 4-8: unchanged least_popular / A,B,C,D rounds (re-affirm values) ---
$ws.Range("B4").Value = "least_popular"
$ws.Range("C4").Value = "A,B,C,D"
$ws.Range("B5").Value = "least_popular"
$ws.Range("C5").Value = "A,B,C,D"
$ws.Range("B6").Value = "least_popular"
$ws.Range("C6").Value = "A,B,C,D"
$ws.Range("B7").Value = "least_popular"
$ws.Range("C7").Value = "A,B,C,D"
$ws.Range("B8").Value = "least_popular"
$ws.Range("C8").Value = "A,B,C,D"

# --- Rows 9-23: new "aggregate_difficulty" rounds with per-row Question numbers and answer letters ---
# Column A values (Question numbers). These cells already carry the text-formatted
# style used by rows 2-8 (style index 4, numFmt "@"), which would make a direct
# Value assignment store the number as text. To keep them genuinely numeric we
# first neutralise the cell format (copy the plain format from A24), write the
# number, then re-apply the real target format (copied from A2) afterwards.
$questionNumbers = @{
    9  = 2.1
    10 = 2.2
    11 = 2.3
    12 = 2.4
    13 = 2.5
    14 = 2.6
    15 = 2.7
    16 = 2.8
    17 = 2.9
    18 = 3.1
    19 = 3.2
    20 = 3.3
    21 = 3.4
    22 = 3.5
    23 = 3.6
}

foreach ($row in 9..23) {
    $ws.Range("A24").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $questionNumbers[$row]
    $ws.Range("A2").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
}

$ws.Range("B9").Value = "aggregate_difficulty"
$ws.Range("C9").Value = "B"

$ws.Range("B10").Value = "aggregate_difficulty"
$ws.Range("C10").Value = "B"

$ws.Range("B11").Value = "aggregate_difficulty"
$ws.Range("C11").Value = "D"

$ws.Range("B12").Value = "aggregate_difficulty"
$ws.Range("C12").Value = "C"

$ws.Range("B13").Value = "aggregate_difficulty"
$ws.Range("C13").Value = "B"

$ws.Range("B14").Value = "aggregate_difficulty"
$ws.Range("C14").Value = "C"

$ws.Range("B15").Value = "aggregate_difficulty"
$ws.Range("C15").Value = "C"

$ws.Range("B16").Value = "aggregate_difficulty"
$ws.Range("C16").Value = "D"

$ws.Range("B17").Value = "aggregate_difficulty"
$ws.Range("C17").Value = "B"

$ws.Range("B18").Value = "aggregate_difficulty"
$ws.Range("C18").Value = "C"

$ws.Range("B19").Value = "aggregate_difficulty"
$ws.Range("C19").Value = "B"

$ws.Range("B20").Value = "aggregate_difficulty"
$ws.Range("C20").Value = "A"

$ws.Range("B21").Value = "aggregate_difficulty"
$ws.Range("C21").Value = "D"

$ws.Range("B22").Value = "aggregate_difficulty"
$ws.Range("C22").Value = "D"

$ws.Range("B23").Value = "aggregate_difficulty"
$ws.Range("C23").Value = "A"

# --- Restore the active-cell selection shown in the saved workbook ---
$ws.Range("C3").Select()
